$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each changed cell value as plain text, preserving the original
# (unstyled) cell formatting: force the Text number format while writing
# so Excel doesn't auto-coerce number-looking strings (e.g. "1.00",
# "612.33") into numeric values, then reset the style back to Normal so
# no stray style index is left on the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.922.77"
Set-TextValue "D3" "3.135.08"
Set-TextValue "E3" "  -3.49%  "
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "612.33"
Set-TextValue "E5" "  +0.36%  "
Set-TextValue "D6" "145.97"
Set-TextValue "E6" "  -7.00%  "
Set-TextValue "E7" "  +0.11%  "
Set-TextValue "D8" "3.131.49"
Set-TextValue "E8" "  -3.60%  "
Set-TextValue "D9" "0.523"
Set-TextValue "E9" "  -3.98%  "
Set-TextValue "E10" "  -6.95%  "
Set-TextValue "D11" "5.29"
Set-TextValue "E11" "  -8.28%  "
Set-TextValue "E12" "  -5.42%  "
Set-TextValue "E13" "  -7.57%  "
Set-TextValue "D14" "35.21"
Set-TextValue "E14" "  -9.68%  "
Set-TextValue "D15" "3.651.29"
Set-TextValue "E15" "  -3.45%  "
Set-TextValue "E16" "  +1.14%  "
Set-TextValue "D17" "63.958.78"
Set-TextValue "E17" "  -3.66%  "
Set-TextValue "D18" "3.139.95"
Set-TextValue "E18" "  -3.55%  "
Set-TextValue "D19" "6.83"
Set-TextValue "E19" "  -8.12%  "
Set-TextValue "D20" "475.27"
Set-TextValue "E20" "  -5.68%  "
Set-TextValue "D21" "14.64"
Set-TextValue "E21" "  -5.11%  "
Set-TextValue "E22" "  -6.77%  "
Set-TextValue "E23" "  -4.10%  "
Set-TextValue "D24" "13.53"
Set-TextValue "E24" "  -7.42%  "
Set-TextValue "D25" "83.25"
Set-TextValue "E25" "  -4.67%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "D27" "2.80"
Set-TextValue "E27" "  -6.90%  "
Set-TextValue "D28" "8.38"
Set-TextValue "E28" "  -8.56%  "
Set-TextValue "D29" "2.15"
Set-TextValue "E29" "  -9.45%  "
Set-TextValue "E30" "  -3.91%  "
Set-TextValue "E31" "  -10.26%  "
Set-TextValue "E32" "  +0.12%  "
Set-TextValue "D33" "2.73"
Set-TextValue "E33" "  -5.26%  "
Set-TextValue "D34" "26.09"
Set-TextValue "E34" "  -6.43%  "
Set-TextValue "D35" "1.12"
Set-TextValue "E35" "  -2.64%  "
Set-TextValue "D36" "5.93"
Set-TextValue "E36" "  -7.97%  "
Set-TextValue "D37" "53.11"
Set-TextValue "E37" "  -4.21%  "
Set-TextValue "E38" "  -6.43%  "
Set-TextValue "D39" "459.74"
Set-TextValue "E39" "  -6.90%  "
Set-TextValue "D40" "2.89"
Set-TextValue "E40" "  -12.42%  "
Set-TextValue "D41" "0.0392"
Set-TextValue "E41" "  -6.74%  "
Set-TextValue "E42" "  -7.87%  "
Set-TextValue "D43" "8.35"
Set-TextValue "E43" "  -5.55%  "
Set-TextValue "D44" "2.835.40"
Set-TextValue "E44" "  -4.95%  "
Set-TextValue "D45" "0.264"
Set-TextValue "E45" "  -9.60%  "
Set-TextValue "E46" "  -10.88%  "
Set-TextValue "E47" "  +0.05%  "
Set-TextValue "B48" "ThetaToken"
Set-TextValue "C48" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D48" "2.36"
Set-TextValue "E48" "  -5.33%  "
Set-TextValue "B49" "InjectiveProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "26.23"
Set-TextValue "E49" "  -9.12%  "
Set-TextValue "E50" "  -4.98%  "
Set-TextValue "D51" "118.22"
Set-TextValue "E51" "  -1.93%  "
